$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.75659999999999
$ws.Range("D3").Value = -7.00599999999999
$ws.Range("B12").Value = 4.866299999999996
$ws.Range("C14").Value = -12.9691
$ws.Range("D20").Value = -7.588399999999999
$ws.Range("D25").Value = -7.464200000000004
$ws.Range("C26").Value = -11.42469999999999
$ws.Range("B27").Value = 5.954500000000007
$ws.Range("D30").Value = -7.063200000000006
$ws.Range("C31").Value = -13.8344
$ws.Range("B32").Value = 6.588799999999999
$ws.Range("C35").Value = -11.67250000000001
$ws.Range("B36").Value = 9.228500000000002
$ws.Range("C37").Value = -12.7933
$ws.Range("B38").Value = 5.0335
$ws.Range("D44").Value = -7.105500000000002
$ws.Range("C45").Value = -13.6834
$ws.Range("B46").Value = 7.012500000000002
$ws.Range("D47").Value = -7.509800000000002
$ws.Range("C52").Value = -10.9898
$ws.Range("B54").Value = 4.6486
$ws.Range("B55").Value = 5.583399999999999
$ws.Range("B56").Value = 4.7411
$ws.Range("C57").Value = -14.47449999999998
$ws.Range("D58").Value = -8.006899999999998
$ws.Range("B67").Value = 5.041899999999993
$ws.Range("B69").Value = 5.086399999999992
$ws.Range("B72").Value = 5.796199999999999
$ws.Range("D78").Value = -7.503900000000003
$ws.Range("C81").Value = -12.90359999999999
$ws.Range("B83").Value = 5.264899999999995
$ws.Range("C83").Value = -13.4948
$ws.Range("D84").Value = -8.692000000000005
$ws.Range("B86").Value = 4.905100000000001
$ws.Range("D89").Value = -6.115799999999997
$ws.Range("B91").Value = 5.0649
$ws.Range("D91").Value = -6.081299999999998
$ws.Range("D92").Value = -5.892699999999999
$ws.Range("B93").Value = 6.767099999999996
$ws.Range("D96").Value = -7.506900000000002
$ws.Range("B99").Value = 4.347299999999998
$ws.Range("C100").Value = -12.41089999999999
$ws.Range("C102").Value = -13.4385
$ws.Range("D102").Value = -8.038599999999997
